$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.126.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.82%  "
$ws.Range("D3").Value = "'2.240.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.75%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'252.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.75%  "
$ws.Range("E6").Value = "  +2.96%  "
$ws.Range("D7").Value = "'75.05"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.76%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "'0.600"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.57%  "
$ws.Range("D10").Value = "'41.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.83%  "
$ws.Range("E11").Value = "  +4.59%  "
$ws.Range("D12").Value = "'6.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.38%  "
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("D14").Value = "'2.576.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'14.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "'2.237.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.07%  "
$ws.Range("D17").Value = "'0.791"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("D18").Value = "'43.025.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.93%  "
$ws.Range("E19").Value = "  +6.03%  "
$ws.Range("D20").Value = "'71.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.67%  "
$ws.Range("D21").Value = "'5.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.24%  "
$ws.Range("D22").Value = "'229.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.67%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "'9.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.62%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "'2.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +15.81%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "'10.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.57%  "
$ws.Range("D27").Value = "'3.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.11%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'2.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.34%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'39.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +28.27%  "
$ws.Range("E30").Value = "  +4.39%  "
$ws.Range("D31").Value = "'171.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("E32").Value = "  +3.87%  "
$ws.Range("E33").Value = "  +7.61%  "
$ws.Range("E34").Value = "  +4.97%  "
$ws.Range("E35").Value = "  +2.35%  "
$ws.Range("E36").Value = "  +10.30%  "
$ws.Range("E37").Value = "  +11.00%  "
$ws.Range("E38").Value = "  +18.68%  "
$ws.Range("E39").Value = "  +11.99%  "
$ws.Range("D40").Value = "'2.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.08%  "
$ws.Range("D41").Value = "'0.205"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.31%  "
$ws.Range("E42").Value = "  +4.47%  "
$ws.Range("D43").Value = "'59.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.68%  "
$ws.Range("E44").Value = "  +6.57%  "
$ws.Range("D45").Value = "'103.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.62%  "
$ws.Range("D46").Value = "'0.480"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +31.14%  "
$ws.Range("D47").Value = "'0.0988"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.31%  "
$ws.Range("D48").Value = "'2.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +15.01%  "
$ws.Range("E49").Value = "  +4.19%  "
$ws.Range("E50").Value = "  +5.20%  "
$ws.Range("E51").Value = "  +3.25%  "
